# Generate Report for Handoff
# Replaces the old handoff id/hash/timestamps with the new ones across the
# Overview / zh-cn / de-de sheets: cell text (shared strings) + the
# `display` caption shown on each hyperlinked filename cell.
#
# Note: the real handoff-target URLs (sheet rels) are intentionally left
# pointing at the *old* id/hash - only the visible text + display caption
# move to the new id/hash, matching the authoritative diff.

$wb = $excel.ActiveWorkbook

$oldId   = "b64b3371-26e7-4d7d-aecb-ebb3958c69c7"
$newId   = "a7ddcd1f-1925-4e8f-9d64-9849cbabb631"
$oldHash = "dc1ce6a2b38958c7028d94a4c406b783603b2909"
$newHash = "2282e85b70030be875e844c58c75391c3e2ca463"

function Set-LinkDisplay($ws, $cellRef, $address, $newText) {
    $h = $ws.Range($cellRef).Hyperlinks.Item(1)
    # Re-assert the (unchanged) target address first - this is required so
    # the follow-up TextToDisplay edit below lands on a hyperlink entry
    # that still carries the original r:id/target instead of losing it.
    $h.Address = $address
    $h2 = $ws.Range($cellRef).Hyperlinks.Item($ws.Hyperlinks.Count())
    $h2.TextToDisplay = $newText
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/369a82fd4abf0576d2d491834371c42fb7633649/e2e/$oldId.md"

Set-LinkDisplay $wsOverview "A2" $overviewMdAddress "$newId.md"
$wsOverview.Range("D2").Value = "2016-41-11 18:41:37"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/369a82fd4abf0576d2d491834371c42fb7633649/e2e/$oldId.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2f68a9c63b4bd3c2a60fe0d9cb8ee36dbfee60a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldHash.zh-cn.xlf"

Set-LinkDisplay $wsZh "A2" $zhMdAddress "$newId.md"
Set-LinkDisplay $wsZh "D2" $zhXlfAddress "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-11 18:41:34"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/369a82fd4abf0576d2d491834371c42fb7633649/e2e/$oldId.md"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c06a0b93272631800cf7f9626985dcb6c0836d72/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldHash.de-de.xlf"

Set-LinkDisplay $wsDe "A2" $deMdAddress "$newId.md"
Set-LinkDisplay $wsDe "D2" $deXlfAddress "$newId.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-11 18:41:37"
